$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D2").Value = "canonical SMILES"

# Data rows: column D gets the "canonical SMILES" (destereoized where applicable)
$ws.Range("D3").Value  = "c1ccc2c(c1)c(ncn2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D4").Value  = "c1ccc2c(c1)c(=Nc3cccc(c3)C(F)(F)F)[nH]cn2"
$ws.Range("D5").Value  = "c1ccc2c(c1)c(=[NH+]c3cccc(c3)C(F)(F)F)nc[nH]2"
$ws.Range("D6").Value  = "c1ccc2c(c1)c([nH+]cn2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D7").Value  = "c1ccc2c(c1)c(ncn2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D8").Value  = "c1ccc2c(c1)c(=Nc3cccc(c3)C(F)(F)F)nc[nH]2"
$ws.Range("D9").Value  = "c1ccc2c(c1)c(nc[nH+]2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D10").Value = "c1ccc2c(c1)c(ncn2)[N-]c3cccc(c3)C(F)(F)F"
$ws.Range("D11").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)Nc3cccc(c3)C(F)(F)F"
$ws.Range("D12").Value = "c1ccc2c(c1)c([nH+]cn2)[NH2+]c3cccc(c3)C(F)(F)F"
$ws.Range("D13").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)[NH2+]c3cccc(c3)C(F)(F)F"

# Copy the formatting (fill/border/font/alignment) from column C onto column D
# for every row so the new column visually matches the rest of the table.
$ws.Range("C2:C13").Copy()
$ws.Range("D2:D13").PasteSpecial(-4122)

# Column width for new column D (target OOXML width 37.7109375;
# Excel's ColumnWidth is quantized to whole pixels internally, so we pick
# the ColumnWidth value that rounds to the closest achievable OOXML width)
$ws.Columns.Item(4).ColumnWidth = 36.875
